$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.39879771986096557,
    -0.0099999988233463455,
    -0.0089999988024338506,
    -0.011999999715694543,
    -0.0059999988172023677,
    -0.0169698125238682,
    -0.019999998584211198,
    -0.01999999857583834,
    -0.0059999987635102059,
    -0.0059999987575949376,
    -0.0044999987783924666,
    -0.0059999987567982416,
    0.0012252853647751749,
    -0.01199999866926138,
    -0.0059999987531362819,
    0.021571763081101203,
    -0.0059999987500942709,
    -0.008999998707504453,
    -0.081933824876212924,
    -0.0089999988241054041,
    -0.0089999988227695837,
    -0.0089999988219533478,
    -0.0089999987925866165,
    -0.041999998315739084,
    -0.041999998307309383,
    -0.0059999987804246757,
    -0.0059999987774870256,
    -0.0059999987631478291,
    -0.011999998672251877,
    -0.019999998557035603,
    0.015964099629115935,
    -0.020999998539487308,
    -0.0059999987483791983
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
